# Update "想去人数" (F column) counts across the four worksheets to match
# the latest scrape output (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 407
$ws1.Range("F5").Value  = 5329
$ws1.Range("F6").Value  = 5329
$ws1.Range("F7").Value  = 169
$ws1.Range("F13").Value = 5280
$ws1.Range("F16").Value = 99
$ws1.Range("F17").Value = 2678
$ws1.Range("F18").Value = 2678
$ws1.Range("F22").Value = 3991
$ws1.Range("F26").Value = 3908
$ws1.Range("F28").Value = 185
$ws1.Range("F29").Value = 252
$ws1.Range("F30").Value = 216
$ws1.Range("F37").Value = 6937
$ws1.Range("F38").Value = 1135
$ws1.Range("F39").Value = 540
$ws1.Range("F40").Value = 973
$ws1.Range("F42").Value = 1411
$ws1.Range("F43").Value = 181
$ws1.Range("F44").Value = 730
$ws1.Range("F45").Value = 25
$ws1.Range("F46").Value = 2364
$ws1.Range("F47").Value = 317
$ws1.Range("F48").Value = 93
$ws1.Range("F49").Value = 12
$ws1.Range("F51").Value = 942

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F25").Value = 829

# --- Sheet "本地生活" (sheet3) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 218

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 218
$ws4.Range("F6").Value  = 407
$ws4.Range("F7").Value  = 5329
$ws4.Range("F8").Value  = 5329
$ws4.Range("F9").Value  = 169
$ws4.Range("F18").Value = 99
$ws4.Range("F19").Value = 2679
$ws4.Range("F20").Value = 2679
$ws4.Range("F24").Value = 3991
$ws4.Range("F25").Value = 3908
$ws4.Range("F27").Value = 185
$ws4.Range("F28").Value = 252
$ws4.Range("F29").Value = 216
$ws4.Range("F36").Value = 6937
$ws4.Range("F37").Value = 1135
$ws4.Range("F38").Value = 540
$ws4.Range("F40").Value = 973
$ws4.Range("F42").Value = 1411
$ws4.Range("F43").Value = 181
$ws4.Range("F44").Value = 730
$ws4.Range("F45").Value = 25
$ws4.Range("F46").Value = 2364
$ws4.Range("F47").Value = 317
$ws4.Range("F48").Value = 93
$ws4.Range("F50").Value = 942
